$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Abril de 2020 a las 23:20"

# Update Cataluña row (row 5) figures
$ws.Range("B5").Value = 24734
$ws.Range("C5").Value = 8635
$ws.Range("D5").Value = 13591
$ws.Range("E5").Value = 2508
